$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.108185291290283
$ws.Range("B1").Value = 4.408207893371582
$ws.Range("C1").Value = 6.987278938293457
$ws.Range("D1").Value = 7.156269550323486
$ws.Range("E1").Value = 5.498712062835693
